# edit.ps1 - apply the "added user and input/output" commit
#
# Summary of the change:
#   1. The footer date field ("24.06.25") used on the slide master and every
#      slide layout is bumped to "27.06.25".
#   2. Three brand-new "Specify System" / "Requirements" slides are added:
#        - a "Users" slide, inserted right after the existing
#          "Intent" slide (new slide #6)
#        - an "Inputs / Outputs" slide, inserted right after that
#          (new slide #7)
#        - a "Requirements" slide, appended at the very end of the deck
#      All previously-existing slides keep their content and simply shift
#      down to make room for the two slides inserted in the middle.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the auto-date footer placeholder everywhere it appears
#    (slide master + every custom layout).
# ---------------------------------------------------------------------
function Update-DateePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Length -gt 0) {
                $tr.Characters(1, $tr.Length).Text = "27.06.25"
            }
        }
    }
}

Update-DateePlaceholder $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateePlaceholder $p.SlideMaster.CustomLayouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# 2. Insert the new slides.
#    NOTE: create the slide that ends up appended at the END of the deck
#    FIRST, then the two that land in the middle - that way the internal
#    slide ids come out in the same order as the target deck (the id
#    counter increments with each new slide, regardless of where it is
#    inserted into the list).
# ---------------------------------------------------------------------
$contentLayout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"

# --- 2a. "Requirements" slide -> appended at the end -------------------
$slideReq = $p.Slides.AddSlide($p.Slides.Count + 1, $contentLayout)
$slideReq.Shapes.Item(1).TextFrame.TextRange.Text = "Requirements"

$bodyReq = $slideReq.Shapes.Item(2).TextFrame.TextRange
$bodyReq.Text = "Functional Requirements:`r`r`r`rNon- Functional Requirements:"
$bodyReq.Paragraphs(4,1).ParagraphFormat.Bullet.Visible = $false

# --- 2b. "Users" slide -> inserted right after slide 5 (Intent) --------
$slideUsers = $p.Slides.AddSlide(6, $contentLayout)
$slideUsers.Shapes.Item(1).TextFrame.TextRange.Text = "Specify System"

$bodyUsers = $slideUsers.Shapes.Item(2).TextFrame.TextRange
$bodyUsers.Text = "Users:`rOn board flight control computer/autopilot: `rexecutes algorithm`rFlight crew:`rReceives stable aircraft and can interact with cockpit control`rMaintenance and health management engineers:`rAnalyse logs`r`r`r"
$bodyUsers.Paragraphs(2,1).IndentLevel = 2
$bodyUsers.Paragraphs(3,1).IndentLevel = 3
$bodyUsers.Paragraphs(4,1).IndentLevel = 2
$bodyUsers.Paragraphs(5,1).IndentLevel = 3
$bodyUsers.Paragraphs(6,1).IndentLevel = 2
$bodyUsers.Paragraphs(7,1).IndentLevel = 3
$bodyUsers.Paragraphs(8,1).IndentLevel = 2
$bodyUsers.Paragraphs(9,1).ParagraphFormat.Bullet.Visible = $false

# --- 2c. "Inputs / Outputs" slide -> inserted right after the "Users" --
$slideIO = $p.Slides.AddSlide(7, $contentLayout)
$slideIO.Shapes.Item(1).TextFrame.TextRange.Text = "Specify System"

$bodyIO = $slideIO.Shapes.Item(2).TextFrame.TextRange
$bodyIO.Text = "Inputs`rAircraft state vector: `rp,q,r,V, alpha, beta, theta, phi, chi, h`rControl surface positions: elevator, rudder, aileron`rGuidance reference: h, theta, phi, beta`rContext data: Fault flags etc.`rOutputs`rInner loop attitude controller: delta(elevator, rudder,aileron)`rOuter loop altitude controller: delta theta`rSupervisory/ status channel: `r`r"
$bodyIO.Paragraphs(1,1).IndentLevel = 2
$bodyIO.Paragraphs(2,1).IndentLevel = 3
$bodyIO.Paragraphs(3,1).IndentLevel = 4
$bodyIO.Paragraphs(4,1).IndentLevel = 3
$bodyIO.Paragraphs(5,1).IndentLevel = 3
$bodyIO.Paragraphs(6,1).IndentLevel = 3
$bodyIO.Paragraphs(7,1).IndentLevel = 2
$bodyIO.Paragraphs(8,1).IndentLevel = 3
$bodyIO.Paragraphs(9,1).IndentLevel = 3
$bodyIO.Paragraphs(10,1).IndentLevel = 3
$bodyIO.Paragraphs(11,1).ParagraphFormat.Bullet.Visible = $false

Write-Host "Slides after edit:" $p.Slides.Count
